$d = $word.ActiveDocument

# 1. Update the text of the "Import the data" bullet run.
$d.Content.Find.Execute(
    " Right click each table, in the order of creation, and import the data from resources/output (make sure Header is selected in Options)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Right click each table, in the order of creation, and import the data from resources/output/csv-files ",
    2)

# 2. Find the paragraph we just edited, then append two new list items after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*csv-files*") {
        $target = $p
        break
    }
}

# --- New sub-bullet: "Make sure Options>Header is selected"
$target.Range.InsertParagraphAfter()
$target = $target.Next()
$target.Range.ListFormat.ListLevelNumber = 2
$target.LeftIndent = 72
$target.FirstLineIndent = -18
$target.Range.Text = "Make sure Options>Header is selected"
$target.Range.Font.Size = 15

# --- New sub-bullet: "Make sure Columns>last_updated is removed"
$target.Range.InsertParagraphAfter()
$target = $target.Next()
$target.Range.ListFormat.ListLevelNumber = 2
$target.LeftIndent = 72
$target.FirstLineIndent = -18
$target.Range.Text = "Make sure Columns>last_updated is removed"
$target.Range.Font.Size = 15
